$d = $word.ActiveDocument

# The name paragraph ("Dheeraj Chand") is the first paragraph in the
# document. The short resume is missing the contact-info line that should
# follow it (before "PROFESSIONAL SUMMARY"). Rebuild paragraph 1 as two
# paragraphs: the existing name paragraph, unchanged, plus a new centered
# paragraph carrying the contact info.
$nameParagraph = $d.Paragraphs(1)
$nameRange = $nameParagraph.Range

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$contact = "202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX"

$xml = "<w:p $ns><w:pPr><w:jc w:val=`"center`"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val=`"28`"/></w:rPr><w:t>Dheeraj Chand</w:t></w:r></w:p>" +
       "<w:p $ns><w:pPr><w:jc w:val=`"center`"/></w:pPr><w:r><w:t>$contact</w:t></w:r></w:p>"

[void]$nameRange.InsertXML($xml)
